# Fixing the big mistake: correct the Total (column B) and Community (column D)
# monthly consumption figures on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = @{ B = 11659.69062328335;  D = 594.91877525 }
    3  = @{ B = 10968.56804391668;  D = 549.4940660333333 }
    4  = @{ B = 11628.53546963335;  D = 567.5385497833333 }
    5  = @{ B = 11296.43561766668;  D = 558.0063456333334 }
    6  = @{ B = 11707.55397813335;  D = 583.4007040333333 }
    7  = @{ B = 11307.89054383335;  D = 567.9765056333333 }
    8  = @{ B = 11706.30873480002;  D = 580.4829874166667 }
    9  = @{ B = 11729.88596511668;  D = 581.9110343833333 }
    10 = @{ B = 11314.81506665002;  D = 563.8412813833334 }
    11 = @{ B = 11694.90437996668;  D = 583.0738891 }
    12 = @{ B = 11334.85010820002;  D = 559.3883812833333 }
    13 = @{ B = 11327.37253641668;  D = 553.8025255333333 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row].B
    $ws.Cells.Item($row, 4).Value = $values[$row].D
}

$wb.Save()
